$d = $word.ActiveDocument

# The document has three <id>...</id> tags whose text is split across
# several runs (e.g. "<id>" / "p052r_1" / "</id>"), each run carrying
# different character formatting. Re-find each split tag and replace it
# with itself so Word COM collapses the matched range into a single run
# that inherits the formatting of the first run in the match (the
# Courier New "<id>"/"</id>" run), exactly as described by the diff.

$ids = @("p052r_1", "p052v_1", "p052v_2")

foreach ($id in $ids) {
    $needle = "<id>" + $id + "</id>"
    $rng = $d.Content
    $found = $rng.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, $needle, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $needle"
    }
}
